# Edit script: insert 3 new "Comercializadora del Agro de Limarí - Tomate" price rows
# at the top of the existing data block (rows 285-287), pushing all subsequent rows
# (285-393) down by three rows (to 288-396). The three rows pushed off the bottom
# (old 391-393) become the new final rows (394-396), i.e. this is just a plain
# row-insert at 285 followed by filling in the values for the 3 brand-new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the block; everything below shifts down by 3.
$ws.Rows("285:287").Insert()

# Constant values shared by every row in this data block.
$colA = 2
$colB = "Comercializadora del Agro de Limarí"
$colC = "Coquimbo"
$colE = 4
$colF = 100112020
$colG = "Tomate"
$colN = "`$/bandeja 18 kilos"
$colO = "Provincia de Limarí"
$colQ = 18
$colR = "Hortaliza"

# Row 285: Larga vida / Primera
$ws.Range("A285").Value = $colA
$ws.Range("B285").Value = $colB
$ws.Range("C285").Value = $colC
$ws.Range("D285").Value = 44510
$ws.Range("E285").Value = $colE
$ws.Range("F285").Value = $colF
$ws.Range("G285").Value = $colG
$ws.Range("H285").Value = "Larga vida"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 1000
$ws.Range("K285").Value = 12500
$ws.Range("L285").Value = 13000
$ws.Range("M285").Value = 12750
$ws.Range("N285").Value = $colN
$ws.Range("O285").Value = $colO
$ws.Range("P285").Value = 708
$ws.Range("Q285").Value = $colQ
$ws.Range("R285").Value = $colR

# Row 286: Larga vida / Segunda
$ws.Range("A286").Value = $colA
$ws.Range("B286").Value = $colB
$ws.Range("C286").Value = $colC
$ws.Range("D286").Value = 44510
$ws.Range("E286").Value = $colE
$ws.Range("F286").Value = $colF
$ws.Range("G286").Value = $colG
$ws.Range("H286").Value = "Larga vida"
$ws.Range("I286").Value = "Segunda"
$ws.Range("J286").Value = 800
$ws.Range("K286").Value = 10500
$ws.Range("L286").Value = 11000
$ws.Range("M286").Value = 10750
$ws.Range("N286").Value = $colN
$ws.Range("O286").Value = $colO
$ws.Range("P286").Value = 597
$ws.Range("Q286").Value = $colQ
$ws.Range("R286").Value = $colR

# Row 287: Larga vida / Tercera
$ws.Range("A287").Value = $colA
$ws.Range("B287").Value = $colB
$ws.Range("C287").Value = $colC
$ws.Range("D287").Value = 44510
$ws.Range("E287").Value = $colE
$ws.Range("F287").Value = $colF
$ws.Range("G287").Value = $colG
$ws.Range("H287").Value = "Larga vida"
$ws.Range("I287").Value = "Tercera"
$ws.Range("J287").Value = 700
$ws.Range("K287").Value = 8500
$ws.Range("L287").Value = 9000
$ws.Range("M287").Value = 8750
$ws.Range("N287").Value = $colN
$ws.Range("O287").Value = $colO
$ws.Range("P287").Value = 486
$ws.Range("Q287").Value = $colQ
$ws.Range("R287").Value = $colR

# Make sure date column keeps date formatting (same number format as other D cells).
$ws.Range("D285:D287").NumberFormat = $ws.Range("D288").NumberFormat
